$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "lightning wk/res" column header to "bolt wk/res"
$ws.Range("N1").Value = "bolt wk/res"

# Add absorb (-1) and reflect (-2) special elemental resistance constants
# for "The Evil Eye" (row 4): null wk/res stays, bolt wk/res (N4) becomes reflect (-2),
# and null wk/res (K4) becomes absorb (-1)
$ws.Range("K4").Value = -1
$ws.Range("N4").Value = -2

# Keep selection on K4 to match the saved workbook view
$ws.Range("K4").Select()
